# ------------------------------------------------------------------
# Adapt createToString and adjustMethodName on PHP type sheet.
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$wsValue = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Insert two new config-summary rows above the inheritance block
#    (old row 10 -> new row 12, everything below shifts down by 2).
# ------------------------------------------------------------------
$wsValue.Rows.Item(10).Resize(2).Insert()

$wsValue.Range("A10").Value = "toStringメソッドの生成"
$wsValue.Range("C10").Value = "○"
$wsValue.Range("A11").Value = "フィールド名の変形"
$wsValue.Range("C11").Value = "○"

Write-Host "done step1"
